# Added a missing label for customer.
# Insert a new row above row 212 on the active sheet and populate it with
# the missing "Kundfordranskonto" (customer receivable account) label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 212 (and everything below it) down by one row.
$ws.Rows.Item(212).Insert()

# Fill in the new row's data.
$ws.Range("A212").Value = "Edit"
$ws.Range("B212").Value = 26598
$ws.Range("D212").Value = "Kundfordranskonto"

# Match the author's final selection (cell D212, where the new label text was entered).
$ws.Range("D212").Select()
